$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows 17-44 (columns B-H) down to rows 18-45
# to make room for a newly-uploaded book entry at row 17.
# Column A (S.No.) is intentionally left untouched: its values don't move,
# only the FileName/Book/Author/Link/Edn/Year/Publisher data shifts.
$src = $ws.Range("B17:H44")
$dst = $ws.Range("B18:H45")
$dst.Value2 = $src.Value2

# Give the newly created row 45's "S.No." cell the same style as the rest
# of column A (bold, bordered, centered) by copying the format from A44.
$ws.Range("A44").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A45").Value2 = 44

# Write the new book entry into row 17.
$ws.Range("B17").Value2 = "Differential Equations and Boundary Value Problems - Computing and Modeling - C. Henry Edwards, David E. Penney, David T. Calvis (2015, Pearson) 5th Edition.pdf"
$ws.Range("C17").Value2 = "Differential Equations and Boundary Value Problems: Computing and Modeling"
$ws.Range("D17").Value2 = "C. Henry Edwards, David E. Penney, David T. Calvis"
$ws.Range("E17").Value2 = "[GitHub](https://raw.githubusercontent.com/malloc42/cool-maths-books/main/Books/Differential%20Equations%20and%20Boundary%20Value%20Problems%20-%20Computing%20and%20Modeling%20-%20C.%20Henry%20Edwards%2C%20David%20E.%20Penney%2C%20David%20T.%20Calvis%20%282015%2C%20Pearson%29%205th%20Edition.pdf)"
$ws.Range("F17").Value2 = "5th"
$ws.Range("G17").Value2 = 2015
$ws.Range("H17").Value2 = "Pearson"

# Clear the copy marquee.
$excel.CutCopyMode = 0
